$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 9 new weekly wellness rows (171-179) -----------------------------
# First tile the formatting (styles/number formats/fonts) of the last
# existing data row (170) down across the new block so dates, fonts and the
# centred "no pain location" look match the rest of the table.
$src = $ws.Range("A170:I170")
$dst = $ws.Range("A171:I179")
$src.Copy($dst)

# Rows 171, 172, 177 and 178 carry a pain-location label in column G, so
# pull that cell's formatting (font) from an existing labelled cell before
# writing the text - otherwise it would keep the "empty/centered" style
# that came along with the row-170 template.
$labelRows = @(171, 172, 177, 178)
foreach ($r in $labelRows) {
    $ws.Range("G160").Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Row data: Date, Player, Volume, Intensite, Fatigue, Douleur, Localisation, Plaisir
$data = @(
    @{ Row=171; Date=45889; Player="Sofiane Belle";  Volume=60; Intensite=6; Fatigue=4; Douleur=2; Loc="Coup tibia ";  Plaisir=6 },
    @{ Row=172; Date=45889; Player="Naim Dhib";       Volume=60; Intensite=5; Fatigue=4; Douleur=3; Loc="Adducteurs "; Plaisir=7 },
    @{ Row=173; Date=45889; Player="Ilyes Boughanmi"; Volume=60; Intensite=6; Fatigue=5; Douleur=0; Loc=$null;         Plaisir=0 },
    @{ Row=174; Date=45889; Player="Omar Benyounes";  Volume=60; Intensite=5; Fatigue=6; Douleur=0; Loc=$null;         Plaisir=6 },
    @{ Row=175; Date=45889; Player="Naim Ighbane";    Volume=60; Intensite=6; Fatigue=6; Douleur=0; Loc=$null;         Plaisir=1 },
    @{ Row=176; Date=45889; Player="Ilan Ihaddadene"; Volume=60; Intensite=7; Fatigue=7; Douleur=0; Loc=$null;         Plaisir=9 },
    @{ Row=177; Date=45889; Player="Emmanuel Valey";  Volume=60; Intensite=7; Fatigue=6; Douleur=8; Loc="Adducteur";   Plaisir=5 },
    @{ Row=178; Date=45889; Player="Karahali Souaré"; Volume=60; Intensite=4; Fatigue=6; Douleur=6; Loc="Cheville";    Plaisir=9 },
    @{ Row=179; Date=45889; Player="Romain Thunet";   Volume=60; Intensite=6; Fatigue=5; Douleur=0; Loc=$null;         Plaisir=7 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 2).Value = $entry.Player
    $ws.Cells.Item($r, 3).Value = $entry.Volume
    $ws.Cells.Item($r, 4).Value = $entry.Intensite
    $ws.Cells.Item($r, 5).Value = $entry.Fatigue
    $ws.Cells.Item($r, 6).Value = $entry.Douleur
    if ($entry.Loc) {
        $ws.Cells.Item($r, 7).Value = $entry.Loc
    }
    $ws.Cells.Item($r, 8).Value = $entry.Plaisir
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"
}

# --- Update the selection / scroll position to reflect the new rows ------
$ws.Range("I170:I179").Select()
